# Add data for 2022-03-07
# - Advance the "through" date in the sheet name and header label from
#   February 26 to February 27, 2022.
# - Update underlying neighborhood-by-month carjacking counts to reflect
#   the newly-added incident(s).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) from "Through 2022-02-26" to "Through 2022-02-27"
$ws.Name = "Through 2022-02-27"

# Update the column header label shared string
$ws.Range("B1").Value = "February 2022 (through February 27)"

# Row 3 - Austin
$ws.Range("D3").Value = 14

# Row 4 - New City
$ws.Range("D4").Value = 2

# Row 11 - Garfield Park
$ws.Range("N11").Value = 4

# Row 15 - West Town
$ws.Range("D15").Value = 6

# Row 17 - Chatham
$ws.Range("F17").Value = 3
$ws.Range("N17").Value = 2

# Row 21 - Logan Square (new value in previously empty cell)
$ws.Range("D21").Value = 1

# Row 27 - Roseland (new value in previously empty cell)
$ws.Range("P27").Value = 1

# Row 29 - Uptown
$ws.Range("D29").Value = 2

# Row 34 - West Loop (new value in previously empty cell)
$ws.Range("D34").Value = 1

# Row 37 - West Pullman
$ws.Range("B37").Value = 3
$ws.Range("D37").Value = 1

# Row 41 - Loop (new value in previously empty cell)
$ws.Range("J41").Value = 1

# Row 50 - Hermosa
$ws.Range("N50").Value = 3

# Row 61 - Chinatown (new value in previously empty cell)
$ws.Range("F61").Value = 1

# Row 68 - Garfield Ridge (new value in previously empty cell)
$ws.Range("D68").Value = 1

# Row 73 - Lincoln Park
$ws.Range("B73").Value = 3

# Row 86 - Washington Park
$ws.Range("B86").Value = 2

# Row 87 - West Elsdon (new value in previously empty cell)
$ws.Range("D87").Value = 1
